$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1441
$ws.Range("I6").Value = 2450
$ws.Range("K6").Value = 7350
$ws.Range("M6").Value = -7238

$ws.Range("H11").Value = 78.111115
$ws.Range("I11").Value = 78.111115
$ws.Range("K11").Value = 78.111115
$ws.Range("M11").Value = 61.888885

$ws.Range("H98").Value = 34021.277
$ws.Range("I98").Value = 38729.24
$ws.Range("K98").Value = 38729.24
$ws.Range("M98").Value = -37231.24

$ws.Range("H116").Value = 6497.7085
$ws.Range("I116").Value = 6359.4546
$ws.Range("J116").Value = 6614.6924
$ws.Range("K116").Value = 6359.4546
$ws.Range("L116").Value = 6614.6924
$ws.Range("M116").Value = -2917.4546
$ws.Range("N116").Value = -13498.6924

$ws.Range("H122").Value = 34021.277
$ws.Range("I122").Value = 38729.24
$ws.Range("K122").Value = 116187.72
$ws.Range("M122").Value = -113737.72

$ws.Range("H138").Value = 2320.4634
$ws.Range("I138").Value = 1568.6428
$ws.Range("J138").Value = 3939.7693
$ws.Range("K138").Value = 4705.928400000001
$ws.Range("L138").Value = 11819.3079
$ws.Range("M138").Value = 434.0715999999993
$ws.Range("N138").Value = -22099.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3052.8262
$ws.Range("I61").Value = 2301.1304
$ws.Range("K61").Value = 2301.1304
$ws.Range("M61").Value = -2089.1304

$ws.Range("H74").Value = 338735.25
$ws.Range("I74").Value = 475172.75
$ws.Range("K74").Value = 475172.75
$ws.Range("M74").Value = -474298.75

$ws.Range("H77").Value = 338735.25
$ws.Range("I77").Value = 475172.75
$ws.Range("K77").Value = 2375863.75
$ws.Range("M77").Value = -2371495.75

$ws.Range("H122").Value = 31082.773
$ws.Range("I122").Value = 1752.1818
$ws.Range("J122").Value = 102779.78
$ws.Range("K122").Value = 5256.5454
$ws.Range("L122").Value = 308339.34
$ws.Range("M122").Value = -2806.5454
$ws.Range("N122").Value = -313239.34

$ws.Range("H136").Value = 3052.8262
$ws.Range("I136").Value = 2301.1304
$ws.Range("K136").Value = 6903.3912
$ws.Range("M136").Value = -4353.3912

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 9999
$ws.Range("I23").Value = 9999
$ws.Range("K23").Value = 9999
$ws.Range("M23").Value = -9759

$ws.Range("H27").Value = 9999
$ws.Range("I27").Value = 9999
$ws.Range("K27").Value = 9999
$ws.Range("M27").Value = -9807

$ws.Range("H31").Value = 1761.2307
$ws.Range("I31").Value = 1665.4062
$ws.Range("K31").Value = 1665.4062
$ws.Range("M31").Value = -1370.4062

$ws.Range("H34").Value = 1761.2307
$ws.Range("I34").Value = 1665.4062
$ws.Range("K34").Value = 1665.4062
$ws.Range("M34").Value = -1463.4062

$ws.Range("H107").Value = 34515364
$ws.Range("I107").Value = 50045716
$ws.Range("K107").Value = 50045716
$ws.Range("M107").Value = -50043796

$ws.Range("H132").Value = 7313.725
$ws.Range("I132").Value = 2703.15
$ws.Range("K132").Value = 8109.450000000001
$ws.Range("M132").Value = -5579.450000000001

$ws.Range("H137").Value = 67999.75
$ws.Range("J137").Value = 67999.75
$ws.Range("L137").Value = 67999.75
$ws.Range("N137").Value = -78199.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 408.46667
$ws.Range("I12").Value = 338.66666
$ws.Range("J12").Value = 455
$ws.Range("K12").Value = 1015.99998
$ws.Range("L12").Value = 1365
$ws.Range("M12").Value = -842.9999799999999
$ws.Range("N12").Value = -1711

$ws.Range("H141").Value = 20030
$ws.Range("I141").Value = 20030
$ws.Range("K141").Value = 60090
$ws.Range("M141").Value = -54910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4916
$ws.Range("I29").Value = 3191
$ws.Range("J29").Value = 6928.5
$ws.Range("K29").Value = 3191
$ws.Range("L29").Value = 6928.5
$ws.Range("M29").Value = -2901
$ws.Range("N29").Value = -7508.5

$ws.Range("H102").Value = 30580.95
$ws.Range("I102").Value = 41391.117
$ws.Range("J102").Value = 10504.929
$ws.Range("K102").Value = 41391.117
$ws.Range("L102").Value = 10504.929
$ws.Range("M102").Value = -39769.117
$ws.Range("N102").Value = -13748.929

$ws.Range("H132").Value = 2941.6775
$ws.Range("I132").Value = 2873.0667
$ws.Range("K132").Value = 8619.2001
$ws.Range("M132").Value = -6089.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1005.7143
$ws.Range("I4").Value = 1010
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 1010
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -897
$ws.Range("N4").Value = -1226

$ws.Range("H5").Value = 23999
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 23999
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 23999
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -24225

$ws.Range("H20").Value = 4905.7144
$ws.Range("J20").Value = 4905.7144
$ws.Range("L20").Value = 4905.7144
$ws.Range("N20").Value = -5357.7144

$ws.Range("H28").Value = 1005.7143
$ws.Range("I28").Value = 1010
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 1010
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -778
$ws.Range("N28").Value = -1464

$ws.Range("H37").Value = 1005.7143
$ws.Range("I37").Value = 1010
$ws.Range("J37").Value = 1000
$ws.Range("K37").Value = 1010
$ws.Range("L37").Value = 1000
$ws.Range("M37").Value = -903
$ws.Range("N37").Value = -1214

$ws.Range("H132").Value = 3080.682
$ws.Range("I132").Value = 2366.7273
$ws.Range("K132").Value = 7100.1819
$ws.Range("M132").Value = -4570.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 26000
$ws.Range("J21").Value = 26000
$ws.Range("L21").Value = 26000
$ws.Range("N21").Value = -26470

$ws.Range("H24").Value = 30995
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 30995
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 30995
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -31455

$ws.Range("H30").Value = 3153.8462

$ws.Range("H35").Value = 26000
$ws.Range("J35").Value = 26000
$ws.Range("L35").Value = 26000
$ws.Range("N35").Value = -26580

$ws.Range("H96").Value = 3625.5
$ws.Range("I96").Value = 3042.8333
$ws.Range("J96").Value = 4499.5
$ws.Range("K96").Value = 3042.8333
$ws.Range("L96").Value = 4499.5
$ws.Range("M96").Value = -1669.8333
$ws.Range("N96").Value = -7245.5

$ws.Range("H126").Value = 35719620
$ws.Range("I126").Value = 41671670
$ws.Range("J126").Value = 7298
$ws.Range("K126").Value = 125015010
$ws.Range("L126").Value = 21894
$ws.Range("M126").Value = -125012540
$ws.Range("N126").Value = -26834

$ws.Range("H132").Value = 4184.84
$ws.Range("I132").Value = 4332
$ws.Range("K132").Value = 12996
$ws.Range("M132").Value = -10466

$ws.Range("H136").Value = 58513.062
$ws.Range("I136").Value = 47515.285
$ws.Range("J136").Value = 77759.164
$ws.Range("K136").Value = 142545.855
$ws.Range("L136").Value = 233277.492
$ws.Range("M136").Value = -139995.855
